# Weekly refresh of the Perejil (parsley) price records.
# For each affected row, update: Fecha (D), Volumen (J), Precio minimo (K),
# Precio maximo (L), Precio promedio ponderado (M), Unidad de comercializacion (N,
# only rows 18/26), Precio $/Kg (P) and Kg o Unidades (Q, only rows 18/26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ Row=2;  D='2020-12-07'; J=200; K=1300; L=1500; M=1400; P=700 },
    @{ Row=3;  D='2021-12-10'; K=900;  L=1000; M=950;  P=475 },
    @{ Row=4;  D='2021-09-29'; J=300; K=900;  L=1000; M=950;  P=475 },
    @{ Row=5;  D='2021-04-05' },
    @{ Row=6;  D='2021-08-19'; J=250; K=1300; L=1500; M=1400; P=700 },
    @{ Row=7;  D='2021-07-26'; K=1800; L=2000; M=1900; P=950 },
    @{ Row=9;  D='2021-09-08'; K=900;  L=1000; M=950;  P=475 },
    @{ Row=10; D='2021-07-13'; J=250; K=2400; L=2500; M=2450; P=1225 },
    @{ Row=11; D='2022-02-09'; J=270; K=2200; L=2500; M=2350; P=1175 },
    @{ Row=12; D='2021-03-11'; K=1700; L=1800; M=1750; P=875 },
    @{ Row=13; D='2020-11-26'; K=900;  L=1000; M=950;  P=475 },
    @{ Row=14; D='2021-02-16'; J=250; K=1200; L=1300; M=1250; P=625 },
    @{ Row=15; D='2022-06-14'; J=250; K=2500; L=2800; M=2650; P=1325 },
    @{ Row=16; D='2021-12-14'; K=900;  L=1000; M=950;  P=475 },
    @{ Row=17; D='2021-07-08'; J=300; K=2400; L=2500; M=2450; P=1225 },
    @{ Row=18; D='2022-01-11'; J=300; K=1400; L=1500; M=1450; N='$/atado 1,5 a 2 kilos'; P=725; Q=2 },
    @{ Row=19; D='2022-02-24'; J=270; K=1300; L=1500; M=1400; P=700 },
    @{ Row=20; D='2021-01-06'; K=1800; L=2000; M=1900; P=950 },
    @{ Row=21; D='2021-08-30'; K=950;  M=975;  P=488 },
    @{ Row=22; D='2021-02-26'; J=250; K=1800; L=2000; M=1900; P=950 },
    @{ Row=23; D='2021-06-16'; J=250; K=2500; L=2800; M=2650; P=1325 },
    @{ Row=25; D='2021-04-16'; J=300; K=900;  L=1000; M=950;  P=475 },
    @{ Row=26; D='2020-12-16'; J=200; K=1000; L=1200; M=1100; N='$/atado'; P=1100; Q=1 },
    @{ Row=27; D='2021-11-25'; J=300; K=1400; L=1500; M=1450; P=725 },
    @{ Row=28; D='2021-06-18'; J=200; K=1800; L=2000; M=1900; P=950 },
    @{ Row=29; D='2021-02-02' },
    @{ Row=30; D='2021-03-02'; J=500; K=1400; L=1500; M=1450; P=725 }
)

foreach ($rec in $rowData) {
    $r = $rec.Row
    foreach ($col in 'D','J','K','L','M','N','P','Q') {
        if ($rec.ContainsKey($col)) {
            $ws.Range("$col$r").Value = $rec[$col]
        }
    }
}
